# Apply BoM reference renumbering + Costs "Created:" timestamp update
$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BoM")

# Row 9: C_Small 68pF -> reference renamed from C5 to C6
$bom.Range("D9").Value = "C6"

# Row 11: C_Small 0.1uF group -> C2 renamed to C5
$bom.Range("D11").Value = "C5 C8 C11 C12 C13"

# Row 12: C_Polarized_Small 10uF group -> C6 renamed to C2
$bom.Range("D12").Value = "C1 C2 C7 C10"

# Row 18: R 100 group -> R2 renamed to R4 (becomes R3 R4)
$bom.Range("D18").Value = "R3 R4"

# Row 19: R 10K group -> R4 renamed to R2 (R2 joins this group)
$bom.Range("D19").Value = "R2 R5 R6 R7 R8 R9 R10 R11"

# Rows 21/22: TL074 (U2) and YAC512 (U1) references swap
$bom.Range("D21").Value = "U1"
$bom.Range("D22").Value = "U2"

$costs = $wb.Worksheets.Item("Costs")

# Update the "Created:" timestamp
$costs.Range("B28").Value = "2025-10-09 21:11:23"
